$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "p1"
$ws.Range("B1").Value = "p2"
$ws.Range("C1").Value = "p3"

$ws.Range("C2").Select()
